# "Generate Report for Handback" -- refresh the localization-status report
# after a de-de/zh-cn handback: the status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the stale-handback warning clears now
# that the handback is current, and the handback timestamps advance.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusNew
$ov.Range("F2").Value = $statusNew

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusNew
$zh.Range("K2").Value = "2016-08-30 16:58:35"
$zh.Range("P2").Value = ""

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(16).ColumnWidth = 12.833333333333334

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusNew
$de.Range("K2").Value = "2016-08-30 16:58:42"
$de.Range("P2").Value = ""

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(16).ColumnWidth = 12.833333333333334

# --- Overview column widths (zh-cn / de-de status columns widened) ---
$ov.Columns.Item(5).ColumnWidth = 29.166666666666668
$ov.Columns.Item(6).ColumnWidth = 29.166666666666668
